# Batch264 new code added
# Append a new "result on<timestamp>" header column (G1) to the register sheet,
# mirroring what the test automation run writes after each execution.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(1, 7).Value = "result on16-Jul-2024-05-35-15"
